$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new item row gets populated ---

# A7: numeric counter 0 -> 1 (style/format unchanged)
$ws.Range("A7").Value = 1

# C7 (merged C7:G7): item name - style s=8 numFmtId changes 0 -> 49 (Text)
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "سرنجات 3 سم"

# H7 (merged H7:K7): style s=9 numFmtId changes 0 -> 49 (Text)
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "-2:0"

# L7 (merged L7:M7): keeps its original custom numeric format (numFmtId 165)
# but the cell content itself becomes the text "0". Force text entry, then
# restore the original number format code so the style stays identical.
$ws.Range("L7:M7").NumberFormat = "@"
$ws.Range("L7").Value = "0"
$ws.Range("L7:M7").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

# N7 (merged N7:O7): shares style s=8 with C7 -> numFmtId 0 -> 49 (Text)
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "2.00"

# Q7: style s=12 numFmtId changes 0 -> 49 (Text)
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# --- Row 8 ---
$ws.Range("P8:Q8").RowHeight = 29.25
$ws.Range("P8").Value = 0

# --- Row 9: refresh footer timestamp (text stays the same style) ---
$ws.Range("A9").Value = "Wednesday, 17 September, 2025 10:29 PM"
